$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A7").Value = 42608.892118055555
$ws.Range("B7").Value = -4
$ws.Range("C7").Value = 56
$ws.Range("D7").Value = 39
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 83
$ws.Range("G7").Value = 26214
$ws.Range("H7").Value = 22394
$ws.Range("I7").Value = 1349
$ws.Range("J7").Value = 211
$ws.Range("K7").Value = 148
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = 15
$ws.Range("N7").Value = "Named"
